$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.615.94'
$ws.Range("E2").Value = '  +2.46%  '
$ws.Range("D3").Value = '2.951.03'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.09'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.56'
$ws.Range("E6").Value = '  +4.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '2.950.24'
$ws.Range("E8").Value = '  +2.24%  '
$ws.Range("E9").Value = '  +2.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.96'
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("E11").Value = '  +8.30%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.435'
$ws.Range("E12").Value = '  +1.40%  '
$ws.Range("E13").Value = '  +6.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.13'
$ws.Range("E14").Value = '  -0.78%  '
$ws.Range("E15").Value = '  -0.87%  '
$ws.Range("D16").Value = '3.440.26'
$ws.Range("E16").Value = '  +2.27%  '
$ws.Range("D17").Value = '62.586.50'
$ws.Range("E17").Value = '  +2.41%  '
$ws.Range("D18").Value = '2.949.78'
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.64'
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '434.00'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.660'
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.94'
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  +5.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.06'
$ws.Range("E25").Value = '  -0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.91'
$ws.Range("E26").Value = '  +4.37%  '
$ws.Range("E27").Value = '  +1.39%  '
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.15'
$ws.Range("E29").Value = '  +5.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.16'
$ws.Range("E30").Value = '  +3.21%  '
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("E32").Value = '  +16.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.108'
$ws.Range("E33").Value = '  +2.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.12'
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.991'
$ws.Range("E36").Value = '  +1.60%  '
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.99'
$ws.Range("E38").Value = '  +5.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '49.68'
$ws.Range("E39").Value = '  +1.09%  '
$ws.Range("E40").Value = '  +4.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.34'
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.274'
$ws.Range("E43").Value = '  +3.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.00'
$ws.Range("E44").Value = '  -4.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '134.81'
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("D46").Value = '2.680.32'
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '352.25'
$ws.Range("E48").Value = '  +1.71%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '22.47'
$ws.Range("E51").Value = '  -1.30%  '
